$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (columns F..J change, 1-indexed 6..10)
$ws.Columns.Item(6).Width = 2.140625
$ws.Columns.Item(7).Width = 3.140625
$ws.Columns.Item(8).Width = 3.140625
$ws.Columns.Item(9).Width = 2.140625
$ws.Columns.Item(10).Width = 3.140625

# Update row 1 values
$ws.Range("C1").Value = 23
$ws.Range("D1").Value = 11
$ws.Range("E1").Value = 30
$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 19
$ws.Range("H1").Value = 32
$ws.Range("I1").Value = 9
$ws.Range("J1").Value = 31
$ws.Range("K1").Value = 10
$ws.Range("M1").Value = 0.033000000000000002
$ws.Range("N1").Value = 0.064000000000000001
$ws.Range("O1").Value = 0.083999999999999991
$ws.Range("P1").Value = 0.088999999999999996
$ws.Range("Q1").Value = 0.085999999999999993
